$wb = $excel.ActiveWorkbook

# --- Sheet "small": add Std_Dev column (G) ---
$small = $wb.Worksheets.Item("small")
$small.Range("G1").Value = "Std_Dev"
$small.Range("G2").Formula = "=STDEV(B2,B5,B8)"
$small.Range("G3").Formula = "=STDEV(B3,B6,B9)"
$small.Range("G4").Formula = "=STDEV(B4,B7,B10)"
[void]$small.Range("G4").Select()

# --- Sheet "large": add Std_Dev column (G) ---
$large = $wb.Worksheets.Item("large")
$large.Range("G1").Value = "Std_Dev"
$large.Range("G2").Formula = "=STDEV(B2,B5,B8)"
$large.Range("G3").Formula = "=STDEV(B3,B6,B9)"
$large.Range("G4").Formula = "=STDEV(B4,B7,B10)"
[void]$large.Range("G4").Select()

# --- Make "experimental_results" the active/selected sheet ---
$results = $wb.Worksheets.Item("experimental_results")
[void]$results.Activate()
